# "Generate Report for Handoff" — localization-status.xlsx update
#
# The two tracked files (30ea3542-...md and fad6f72e-...md) swap rows
# (row 2 <-> row 3) on every sheet because fad6f72e is now the file that
# was handed back earlier and 30ea3542 becomes "Ready for handoff" with a
# fresh timestamp and an out-of-date-handback error message.

$wb = $excel.ActiveWorkbook

function Set-LinkDisplay($ws, $addr, $text) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 now describes fad6f72e (previously described 30ea3542)
$ov.Range("A2").Value = "fad6f72e-85ae-4231-a2be-c918fcf075b8.md"
$ov.Range("B2").Value = "e2e\fad6f72e-85ae-4231-a2be-c918fcf075b8.md"
Set-LinkDisplay $ov '$B$2' "e2e\fad6f72e-85ae-4231-a2be-c918fcf075b8.md"

# Row 3 now describes 30ea3542, which just moved to "Ready for handoff"
$ov.Range("A3").Value = "30ea3542-e5e4-4061-8def-87899eb27948.md"
$ov.Range("B3").Value = "e2e\30ea3542-e5e4-4061-8def-87899eb27948.md"
Set-LinkDisplay $ov '$B$3' "e2e\30ea3542-e5e4-4061-8def-87899eb27948.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-03 22:51:52"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "fad6f72e-85ae-4231-a2be-c918fcf075b8.md"
$zh.Range("G2").Value = "fad6f72e-85ae-4231-a2be-c918fcf075b8.7ee745c1015371e206be0df10ee423386e6ace5b.zh-cn.xlf"
$zh.Range("I2").Value = "fad6f72e-85ae-4231-a2be-c918fcf075b8.md"
$zh.Range("J2").Value = "fad6f72e-85ae-4231-a2be-c918fcf075b8.7ee745c1015371e206be0df10ee423386e6ace5b.zh-cn.xlf"
Set-LinkDisplay $zh '$A$2' "fad6f72e-85ae-4231-a2be-c918fcf075b8.md"
Set-LinkDisplay $zh '$I$2' "fad6f72e-85ae-4231-a2be-c918fcf075b8.md"

$zh.Range("A3").Value = "30ea3542-e5e4-4061-8def-87899eb27948.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "30ea3542-e5e4-4061-8def-87899eb27948.c74d6c3641e984933863352a5f02be88d5fb1ff0.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-03 22:51:48"
$zh.Range("I3").Value = "30ea3542-e5e4-4061-8def-87899eb27948.md"
$zh.Range("J3").Value = "30ea3542-e5e4-4061-8def-87899eb27948.c74d6c3641e984933863352a5f02be88d5fb1ff0.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/79f9d5588b399f008fac7d885a89e29135e410e6/e2e/30ea3542-e5e4-4061-8def-87899eb27948.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef6224dfc08e0a716073ab1f8a52de282e9b0ef9/e2e/30ea3542-e5e4-4061-8def-87899eb27948.md."
Set-LinkDisplay $zh '$A$3' "30ea3542-e5e4-4061-8def-87899eb27948.md"
Set-LinkDisplay $zh '$I$3' "30ea3542-e5e4-4061-8def-87899eb27948.md"

$zh.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "fad6f72e-85ae-4231-a2be-c918fcf075b8.md"
$de.Range("G2").Value = "fad6f72e-85ae-4231-a2be-c918fcf075b8.7ee745c1015371e206be0df10ee423386e6ace5b.de-de.xlf"
$de.Range("I2").Value = "fad6f72e-85ae-4231-a2be-c918fcf075b8.md"
$de.Range("J2").Value = "fad6f72e-85ae-4231-a2be-c918fcf075b8.7ee745c1015371e206be0df10ee423386e6ace5b.de-de.xlf"
Set-LinkDisplay $de '$A$2' "fad6f72e-85ae-4231-a2be-c918fcf075b8.md"
Set-LinkDisplay $de '$I$2' "fad6f72e-85ae-4231-a2be-c918fcf075b8.md"

$de.Range("A3").Value = "30ea3542-e5e4-4061-8def-87899eb27948.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "30ea3542-e5e4-4061-8def-87899eb27948.c74d6c3641e984933863352a5f02be88d5fb1ff0.de-de.xlf"
$de.Range("H3").Value = "2016-09-03 22:51:52"
$de.Range("I3").Value = "30ea3542-e5e4-4061-8def-87899eb27948.md"
$de.Range("J3").Value = "30ea3542-e5e4-4061-8def-87899eb27948.c74d6c3641e984933863352a5f02be88d5fb1ff0.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/79f9d5588b399f008fac7d885a89e29135e410e6/e2e/30ea3542-e5e4-4061-8def-87899eb27948.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef6224dfc08e0a716073ab1f8a52de282e9b0ef9/e2e/30ea3542-e5e4-4061-8def-87899eb27948.md."
Set-LinkDisplay $de '$A$3' "30ea3542-e5e4-4061-8def-87899eb27948.md"
Set-LinkDisplay $de '$I$3' "30ea3542-e5e4-4061-8def-87899eb27948.md"

$de.Columns.Item(16).ColumnWidth = 39.17
